$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new column value (B4), which also extends the shared-string
# table with "another change" and grows the used range to A1:B4.
$ws.Range("B4").Value = "another change"

# Move/refresh the sheet's active selection to the newly added cell.
$null = $ws.Range("B4").Select()
